# New crime data collected — weekly refresh of the 42nd Precinct CompStat
# report: bumps the report "Volume/Number" and the covered week dates in
# the header, then overwrites the weekly/28-day/YTD/2-year crime-complaint
# statistics table (rows 15-30) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 31   Number  49" -> "Volume 31   Number  50"
# and     "Report Covering the Week  12/2/2024  Through  12/8/2024"
#     ->  "Report Covering the Week  12/9/2024  Through  12/15/2024"
# These are rich-text shared strings; COM Value assignment collapses the
# run formatting, but reproduces the visible text faithfully.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/9/2024  Through  12/15/2024"

# ---------------------------------------------------------------------
# Helper: write a text value into a cell while preserving the "General"
# number-format style used by its sibling text cells (e.g. the ones that
# already show "0" / "***.*" for N/A rows). We force Text storage via the
# "@" number format so numeric-looking strings like "0" aren't coerced
# back into numbers, then paste-special the formats from a neighboring
# cell that already carries the correct style.
# ---------------------------------------------------------------------
function Set-TextCell($addr, $text, $styleSourceAddr) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Row 15 - Rape
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 40
$ws.Range("K15").Value = 29.032258064516
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 90.476190476190
$ws.Range("N15").Value = -45.205479452054

# Row 16 - Robbery
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 46
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 53.333333333333
$ws.Range("I16").Value = 450
$ws.Range("J16").Value = 408
$ws.Range("K16").Value = 10.294117647058
$ws.Range("L16").Value = 6.635071090047
$ws.Range("M16").Value = 45.631067961165
$ws.Range("N16").Value = -63.942307692307

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 112.5
$ws.Range("F17").Value = 51
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 707
$ws.Range("J17").Value = 678
$ws.Range("K17").Value = 4.277286135693
$ws.Range("L17").Value = 8.103975535168
$ws.Range("M17").Value = 129.545454545455
$ws.Range("N17").Value = -24.304068522483

# Row 18 - Burglary
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -5.263157894736
$ws.Range("I18").Value = 278
$ws.Range("J18").Value = 198
$ws.Range("K18").Value = 40.404040404040
$ws.Range("L18").Value = -7.023411371237
$ws.Range("M18").Value = 89.115646258503
$ws.Range("N18").Value = -74.042950513538

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 528
$ws.Range("J19").Value = 501
$ws.Range("K19").Value = 5.389221556886
$ws.Range("L19").Value = 5.6
$ws.Range("M19").Value = 132.599118942731
$ws.Range("N19").Value = 47.899159663865

# Row 20 - G.L.A.
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -7.142857142857
$ws.Range("I20").Value = 257
$ws.Range("J20").Value = 405
$ws.Range("K20").Value = -36.543209876543
$ws.Range("L20").Value = -7.885304659498
$ws.Range("M20").Value = 125.438596491228
$ws.Range("N20").Value = -43.763676148796

# Row 21 - TOTAL
$ws.Range("C21").Value = 47
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 46.875
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 159
$ws.Range("H21").Value = 10.062893081761
$ws.Range("I21").Value = 2265
$ws.Range("J21").Value = 2233
$ws.Range("K21").Value = 1.433049708911
$ws.Range("L21").Value = 2.907769195820
$ws.Range("M21").Value = 99.033391915641
$ws.Range("N21").Value = -45.709491850431

# Row 22 - Transit (C, G, H switch from numbers to the "N/A" text markers
# used elsewhere in the sheet: "0" and "***.*")
Set-TextCell "C22" "0" "D22"
Set-TextCell "G22" "0" "D22"
Set-TextCell "H22" "***.*" "E22"
$ws.Range("M22").Value = -35.294117647058

# Row 23 - Housing
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 25
$ws.Range("F23").Value = 21
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = -8.695652173913
$ws.Range("I23").Value = 377
$ws.Range("J23").Value = 404
$ws.Range("K23").Value = -6.683168316831
$ws.Range("L23").Value = 7.102272727272
$ws.Range("M23").Value = 96.354166666666

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 45
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 30.487804878048
$ws.Range("I24").Value = 1020
$ws.Range("J24").Value = 1138
$ws.Range("K24").Value = -10.369068541300
$ws.Range("L24").Value = -17.274939172749
$ws.Range("M24").Value = 36

# Row 25 - Retail Theft
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 200
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = -7.142857142857
$ws.Range("I25").Value = 165
$ws.Range("J25").Value = 255
$ws.Range("K25").Value = -35.294117647058
$ws.Range("L25").Value = -60.714285714285

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 81
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = 35
$ws.Range("I26").Value = 962
$ws.Range("J26").Value = 998
$ws.Range("K26").Value = -3.607214428857
$ws.Range("L26").Value = 3.218884120171
$ws.Range("M26").Value = 13.043478260869

# Row 27 - UCR Rape*
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 56
$ws.Range("K27").Value = 7.692307692307
$ws.Range("L27").Value = -11.111111111111

# Row 28 - Other Sex Crimes
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("J28").Value = 92
$ws.Range("K28").Value = -18.478260869565

# Row 29 - Shooting Vic. (only the 2-Year % Chg moves)
$ws.Range("N29").Value = -62.393162393162

# Row 30 - Shooting Inc. (only the 2-Year % Chg moves)
$ws.Range("N30").Value = -65.178571428571
